$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("2025-09-18")
$src.Copy($null, $src)
$ws = $wb.Worksheets.Item($src.Index + 1)
$ws.Name = "2025-09-19"
$wb.Worksheets.Item(1).Activate()

$data = @(
    @("ワンパンマン", "原作/ＯＮＥ 作画/村田雄介", "211撃目"),
    @("追放されたチート付与魔術師は 気ままなセカンドライフを謳歌する。", "六志麻あさ 業務用餅 kisui", "第７１話"),
    @("落ちこぼれだった兄が実は最強 ～史上最強の勇者は転生し、学園で無自覚に無双する～", "村上よしゆき 茨木野 あるてら", "第４２話　勇者、六邪神将相手に舐めプしてたら、ピンチになる（３）"),
    @("ライドンキング", "馬場康誌", "第83話 大統領と龍の闇卵（前編）"),
    @("俺以外誰も採取できない素材なのに「素材採取率が低い」とパワハラする幼馴染錬金術師と絶縁した専属魔導士、辺境の町でスローライフを送りたい。", "狐御前(原作) 西岡知三(作画) ＮＯＣＯ(キャラクター原案)", "第25話-2"),
    @("最弱貴族に転生したので悪役たちを集めてみた", "空野進 sorani ファルまろ", "第12話　最弱貴族、悪役令嬢と出会う（３）"),
    @("美人女上司滝沢さん", "やんBARU(著者)", "第203話"),
    @("黄金の経験値", "原純(原作) 霜月汐(作画) fixro2n(キャラクター原案)", "第18話（前編）"),
    @("デスゲームに巻き込まれた山本さん、気ままにゲームバランスを崩壊させる", "ぽち(原作) カモトタツヤ(作画) 久賀フーナ(キャラクター原案)", "第6話（後編）"),
    @("35歳独身山田、異世界村に理想のセカンドハウスを作りたい　～異世界と現実のいいとこどりライフ～", "出雲大吉(原作) 西尾洋一(作画) ゆのひと(キャラクター原案)", "第4話①"),
    @("時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―", "光永康則", "第６９話『岩鬼停止』③"),
    @("地元のいじめっ子達に仕返ししようとしたら、別の戦いが始まった。", "マツモトケンゴ", "第64話 更衣室の戦いが始まった（２）"),
    @("実は俺、最強でした？", "原作：澄守 彩 漫画：高橋 愛", "第126話　王家の立場"),
    @("蜘蛛ですが、なにか？", "かかし朝浩(著者) 馬場翁(原作) 輝竜司(キャラクター原案)", "第76話その1"),
    @("仮面の黒騎士。正体バレたのでもう学園でも無双する", "楓原こうた(原作) さみ(作画) へいろー(キャラクター原案)", "第2話 王立学園"),
    @("極振り拒否して手探りスタート！　特化しないヒーラー、仲間と別れて旅に出る", "蒼井一秀(作画) 刻一(原作) MIYA*KI(キャラクター原案)", "第66話"),
    @("ゲーム悪役貴族に転生した俺は、チート筋肉で無双する", "昼行燈（原作） しいたけ元帥（漫画）", "特別イラスト"),
    @("双子まとめて『カノジョ』にしない？", "飴色みそ(漫画) 白井ムク(原作) 千種みのり(キャラクター原案)", "ツイントーク：第5話"),
    @("となりの席のヤツがそういう目で見てくる", "mmk", "第50話 ベッドの上で"),
    @("貧乏騎士に嫁入りしたはずが!? 〜野人令嬢は皇太子妃になっても竜を狩りたい〜", "漫画：夏川そぞろ 原作：宮前葵 キャラクター原案：ののまろ", "第12話③立太子式"),
    @("転生大聖女の異世界のんびり紀行", "キダニエル 四葉夕ト", "第５０話　ホリーの“聖槍”です！"),
    @("『おっさんは荷物でも持ってろよ』と新人に舐められてるけど、実は最強の剣豪です", "kiri 日之影ソラ", "第3話 後編"),
    @("規格外のダンジョン攻略者、実は異世界帰りの元勇者", "作画：やまざき君 原作：榊与一", "第7話(2)"),
    @("元・世界１位のサブキャラ育成日記 ～廃プレイヤー、異世界を攻略中！～", "沢村治太郎(原作) 前田理想(漫画) まろ(キャラクター原案)", "第78話その2"),
    @("戦姫サバイバルサガ-異世界の運命をかけた無人島フジュン異性交遊-", "OTOSAMA(著者)", "第22話"),
    @("異世界魔王と召喚少女の奴隷魔術", "原作：むらさきゆきや 漫画：福田直叶 キャラクター原案：鶴崎貴大", "第128話　レムと話してみる（後編）"),
    @("ルパン三世 異世界の姫君（ネイバーワールドプリンセス）", "モンキー・パンチ／エム・ピー・ワークス 内々けやき 佐伯庸介 白狼", "第110話：王女に贈る子守歌"),
    @("やり込んだ乙女ゲームの悪役モブですが、断罪は嫌なので真っ当に生きます@COMIC", "戸張ちょも（漫画） MIZUNA（原作） Ruki（キャラクター原案）", "第20話"),
    @("僕のカノジョ先生", "星河蟹(作画) 孟倫（ＳＤｗｉｎｇ）(構成) 鏡遊(原作) おりょう(キャラクター原案)", "77時間目-1"),
    @("異世界転移で女神様から祝福を！～いえ、手持ちの異能があるので結構です～@COMIC", "コーダ 壁アキオ", "第5話「連携戦闘③」"),
    @("独身貴族は異世界を謳歌する ～結婚しない男の優雅なおひとりさまライフ～", "漫画：駒鳥 ひわ 原作：錬金王 キャラクター原案：三登 いつき", "第34話 独身貴族は礼の品を贈る（1）"),
    @("かつての暗殺者は来世で違う生き方をする", "ツネ(漫画) 丘野優(原作) つなかわ(キャラクター原案)", "第6話②"),
    @("聖者無双", "漫画：秋風緋色 原作：ブロッコリーライオン キャラクター原案：sime", "第92話　龍と精霊の信仰（後半）"),
    @("生徒会役員共", "氏家ト全", "#414"),
    @("転生コロシアム～最弱スキルで最強の女たちを攻略して奴隷ハーレム作ります～", "zunta(作画) はらわたさいぞう(原作)", "第32話：思考を奪う③"),
    @("直径3cmの召喚陣<リミットリング>で「雑魚すら呼べない」と蔑まれた底辺召喚士が頂点に立つまで", "作画：まっつー 原作：空松蓮司", "第7話(2)"),
    @("役目を果たした日陰の勇者は、辺境で自由に生きていきます", "船野真帆(作画) 丘野優(原作) 布施龍太(キャラクター原案)", "第7話前半"),
    @("勇者パーティーをクビになったので故郷に帰ったら、メンバー全員がついてきたんだが", "絶叫あいす。(漫画) 木の芽(原作) 希(キャラクター原案)", "第5話 前編"),
    @("宇崎ちゃんは遊びたい！", "丈(著者)", "第127話"),
    @("魔法少女リリカルなのは EXCEEDS", "都築真紀 川上修一", "番外編②"),
    @("ゲーム内最強の『裏ボス』に転生したので、主人公の代わりに最速クリアを目指します！", "作画：こめぐ 原作：迅空也", "第7話(2)"),
    @("ひげを剃る。そして女子高生を拾う。", "足立いまる(漫画) しめさば(原作) ぶーた(キャラクター原案)", "最終話その1"),
    @("他校の氷姫を助けたら、お友達から始める事になりました", "椀田くろ(作画) 皐月 陽龍(原作) みすみ(キャラクター原案)", "第11話-2"),
    @("オタクに優しいギャルはいない!?", "のりしろちゃん 魚住さかな", "【#157】くらいつけ！ 〜レッツバイト〜"),
    @("カナン様はあくまでチョロい", "nonco", "第151話`tカナンと供犠の印"),
    @("新米オッサン冒険者、最強パーティに死ぬほど鍛えられて無敵になる", "漫画：荻野ケン 原作：岸馬きらく キャラクター原案：Tea", "第71話"),
    @("異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～", "夜桜ユノ 佐野綾佑", "第3話 存在価値"),
    @("勇者は魔王が好きらしい", "赤槻コウ(著者)", "第6話-2：囚われの勇者"),
    @("辺境モブ貴族のウチに嫁いできた悪役令嬢が、めちゃくちゃできる良い嫁なんだが？", "tera(原作) 朝倉はやて(作画) 徹田(キャラクター原案)", "番外編①"),
    @("よわよわ先生", "福地カミオ", "第116話`tほどほどのホワイトデー")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $i + 1
    $ws.Cells.Item($r, 2).Value = $data[$i][0]
    $ws.Cells.Item($r, 3).Value = $data[$i][1]
    $ws.Cells.Item($r, 4).Value = $data[$i][2]
}